$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.020.22'
$ws.Range("E2").Value = '  -5.98%  '
$ws.Range("D3").Value = '2.448.46'
$ws.Range("E3").Value = '  -8.48%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = "'539.60"
$ws.Range("E5").Value = '  -2.52%  '
$ws.Range("D6").Value = "'147.24"
$ws.Range("E6").Value = '  -6.77%  '
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = '  -0.22%  '
$ws.Range("E8").Value = '  -3.40%  '
$ws.Range("D9").Value = '2.465.26'
$ws.Range("E9").Value = '  -7.93%  '
$ws.Range("E10").Value = '  -6.19%  '
$ws.Range("E11").Value = '  -2.02%  '
$ws.Range("D12").Value = "'5.34"
$ws.Range("E12").Value = '  -0.36%  '
$ws.Range("E13").Value = '  -4.61%  '
$ws.Range("D14").Value = '2.887.23'
$ws.Range("E14").Value = '  -8.37%  '
$ws.Range("D15").Value = "'23.99"
$ws.Range("E15").Value = '  -9.27%  '
$ws.Range("D16").Value = '58.917.20'
$ws.Range("E16").Value = '  -6.03%  '
$ws.Range("E17").Value = '  -6.18%  '
$ws.Range("D18").Value = '2.515.52'
$ws.Range("E18").Value = '  -6.18%  '
$ws.Range("D19").Value = "'11.10"
$ws.Range("E19").Value = '  -6.44%  '
$ws.Range("E20").Value = '  -5.94%  '
$ws.Range("D21").Value = "'324.90"
$ws.Range("E21").Value = '  -5.74%  '
$ws.Range("D22").Value = "'0.967"
$ws.Range("E22").Value = '  -3.20%  '
$ws.Range("E23").Value = '  -9.14%  '
$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").Value = "'60.71"
$ws.Range("E24").Value = '  -4.12%  '
$ws.Range("B25").Value = 'Polygon'
$ws.Range("C25").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D25").Value = "'0.455"
$ws.Range("E25").Value = '  -10.24%  '
$ws.Range("E26").Value = '  -4.93%  '
$ws.Range("E27").Value = '  -2.16%  '
$ws.Range("D28").Value = "'7.70"
$ws.Range("E28").Value = '  -6.38%  '
$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D29").Value = "'1.82"
$ws.Range("E29").Value = '  -5.99%  '
$ws.Range("B30").Value = 'PEPE'
$ws.Range("C30").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D30").Value = '0.0₃0772'
$ws.Range("E30").Value = '  -9.42%  '
$ws.Range("B31").Value = 'Aptos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D31").Value = "'6.69"
$ws.Range("E31").Value = '  -7.56%  '
$ws.Range("B32").Value = 'Fetch.AI'
$ws.Range("C32").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D32").Value = "'1.26"
$ws.Range("E32").Value = '  -11.88%  '
$ws.Range("E33").Value = '  -0.11%  '
$ws.Range("D34").Value = "'156.61"
$ws.Range("E34").Value = '  -4.06%  '
$ws.Range("E35").Value = '  -7.73%  '
$ws.Range("D36").Value = "'18.43"
$ws.Range("E36").Value = '  -5.34%  '
$ws.Range("D37").Value = "'4.46"
$ws.Range("E37").Value = '  -9.03%  '
$ws.Range("E38").Value = '  -3.87%  '
$ws.Range("D39").Value = "'317.12"
$ws.Range("E39").Value = '  -10.22%  '
$ws.Range("D40").Value = "'5.83"
$ws.Range("E40").Value = '  -5.93%  '
$ws.Range("B41").Value = 'SuiNetwork'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D41").Value = "'0.837"
$ws.Range("E41").Value = '  -11.68%  '
$ws.Range("B42").Value = 'OKB'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D42").Value = "'36.21"
$ws.Range("E42").Value = '  -5.66%  '
$ws.Range("E43").Value = '  -7.04%  '
$ws.Range("D44").Value = "'0.997"
$ws.Range("E44").Value = '  -0.23%  '
$ws.Range("D45").Value = "'10.73"
$ws.Range("E45").Value = '  -2.65%  '
$ws.Range("E46").Value = '  -3.06%  '
$ws.Range("E47").Value = '  -5.98%  '
$ws.Range("E48").Value = '  -6.28%  '
$ws.Range("E49").Value = '  -5.44%  '
$ws.Range("D50").Value = "'121.66"
$ws.Range("E50").Value = '  -5.44%  '
$ws.Range("D51").Value = "'18.85"
$ws.Range("E51").Value = '  -9.80%  '
